function Set-TextCell {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextCell $ws 2 4 '22.404.42'
$ws.Cells.Item(2, 5).Value = '  -4.68%  '

Set-TextCell $ws 3 4 '1.572.36'
$ws.Cells.Item(3, 5).Value = '  -4.64%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 5).Value = '  +0.06%  '

Set-TextCell $ws 6 4 '291.46'
$ws.Cells.Item(6, 5).Value = '  -2.86%  '

Set-TextCell $ws 7 4 '0.3677'
$ws.Cells.Item(7, 5).Value = '  -3.22%  '

Set-TextCell $ws 8 4 '49.54'
$ws.Cells.Item(8, 5).Value = '  -2.31%  '

Set-TextCell $ws 9 4 '0.3371'
$ws.Cells.Item(9, 5).Value = '  -5.60%  '

$ws.Cells.Item(10, 5).Value = '  -4.68%  '

Set-TextCell $ws 11 4 '0.07576'
$ws.Cells.Item(11, 5).Value = '  -6.62%  '

Set-TextCell $ws 12 4 '1.001'
$ws.Cells.Item(12, 5).Value = '  -0.02%  '

Set-TextCell $ws 13 4 '21.11'
$ws.Cells.Item(13, 5).Value = '  -4.56%  '

Set-TextCell $ws 14 4 '6.048'
$ws.Cells.Item(14, 5).Value = '  -5.85%  '

Set-TextCell $ws 15 4 '6.858'
$ws.Cells.Item(15, 5).Value = '  -7.60%  '

$ws.Cells.Item(16, 2).Value = 'ShibaInu'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 16 4 '0.00001139'
$ws.Cells.Item(16, 5).Value = '  -5.49%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws 17 4 '1.575.33'
$ws.Cells.Item(17, 5).Value = '  -4.84%  '

Set-TextCell $ws 18 4 '89.40'

Set-TextCell $ws 19 4 '0.06753'
$ws.Cells.Item(19, 5).Value = '  -3.35%  '

$ws.Cells.Item(20, 5).Value = '  +0.14%  '

Set-TextCell $ws 21 4 '6.230'
$ws.Cells.Item(21, 5).Value = '  -8.05%  '

$ws.Cells.Item(22, 5).Value = '  -6.57%  '

Set-TextCell $ws 23 4 '11.94'
$ws.Cells.Item(23, 5).Value = '  -5.09%  '

$ws.Cells.Item(24, 2).Value = 'WrappedBTC'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell $ws 24 4 '22.412.29'
$ws.Cells.Item(24, 5).Value = '  -4.71%  '

$ws.Cells.Item(25, 2).Value = 'Toncoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 25 4 '2.415'
$ws.Cells.Item(25, 5).Value = '  -2.68%  '

Set-TextCell $ws 26 4 '2.964'
$ws.Cells.Item(26, 5).Value = '  +1.63%  '

Set-TextCell $ws 27 4 '19.78'
$ws.Cells.Item(27, 5).Value = '  -5.76%  '

Set-TextCell $ws 28 4 '145.87'
$ws.Cells.Item(28, 5).Value = '  -4.44%  '

Set-TextCell $ws 29 4 '4.926'
$ws.Cells.Item(29, 5).Value = '  -5.79%  '

Set-TextCell $ws 30 4 '125.12'
$ws.Cells.Item(30, 5).Value = '  -6.18%  '

Set-TextCell $ws 31 4 '1.748.97'
$ws.Cells.Item(31, 5).Value = '  -4.87%  '

Set-TextCell $ws 32 4 '6.261'
$ws.Cells.Item(32, 5).Value = '  -9.88%  '

Set-TextCell $ws 33 4 '1.982'
$ws.Cells.Item(33, 5).Value = '  -7.49%  '

Set-TextCell $ws 34 4 '0.9874'
$ws.Cells.Item(34, 5).Value = '  -4.32%  '

Set-TextCell $ws 35 4 '10.44'
$ws.Cells.Item(35, 5).Value = '  -12.77%  '

Set-TextCell $ws 36 4 '0.08452'
$ws.Cells.Item(36, 5).Value = '  -3.39%  '

Set-TextCell $ws 37 4 '0.02542'
$ws.Cells.Item(37, 5).Value = '  -6.97%  '

Set-TextCell $ws 38 4 '0.2297'
$ws.Cells.Item(38, 5).Value = '  -6.48%  '

Set-TextCell $ws 39 4 '0.06525'
$ws.Cells.Item(39, 5).Value = '  -5.10%  '

Set-TextCell $ws 40 4 '5.500'
$ws.Cells.Item(40, 5).Value = '  -8.08%  '

Set-TextCell $ws 41 4 '1.261'
$ws.Cells.Item(41, 5).Value = '  -4.80%  '

$ws.Cells.Item(42, 5).Value = '  -12.11%  '

Set-TextCell $ws 43 4 '0.6380'
$ws.Cells.Item(43, 5).Value = '  -7.89%  '

Set-TextCell $ws 44 4 '14.36'
$ws.Cells.Item(44, 5).Value = '  -8.81%  '

$ws.Cells.Item(45, 5).Value = '  +0.05%  '

Set-TextCell $ws 46 4 '0.5992'
$ws.Cells.Item(46, 5).Value = '  -7.16%  '

Set-TextCell $ws 47 4 '3.774'
$ws.Cells.Item(47, 5).Value = '  -3.93%  '

Set-TextCell $ws 48 4 '2.115'
$ws.Cells.Item(48, 5).Value = '  -7.01%  '

Set-TextCell $ws 49 4 '121.04'
$ws.Cells.Item(49, 5).Value = '  -5.69%  '

$ws.Cells.Item(50, 2).Value = 'EOS'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws 50 4 '1.194'
$ws.Cells.Item(50, 5).Value = '  +0.34%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 51 4 '0.07282'
$ws.Cells.Item(51, 5).Value = '  -6.90%  '

Write-Output "Update complete"